$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# Sheet "Overview" (sheet1): add two rows (File Name, Path And Name, ...)
# -----------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$tOverview = $wsOverview.ListObjects.Item("Overview")
$tOverview.ListRows.Add() | Out-Null
$tOverview.ListRows.Add() | Out-Null

# Row 4
$wsOverview.Range("A4").Value = "787dd08a-b847-4b32-8a73-c840274ea705.md"
$wsOverview.Range("B4").Value = "e2e\787dd08a-b847-4b32-8a73-c840274ea705.md"
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("D4").Value = "'"
$wsOverview.Range("E4").Value = "Ready for handoff"
$wsOverview.Range("F4").Value = "Ready for handoff"
$wsOverview.Range("G4").Value = "2016-09-01 06:46:14"

# Row 5
$wsOverview.Range("A5").Value = "c11397f3-0a73-4974-93fe-f54f7fc4bc9d.md"
$wsOverview.Range("B5").Value = "e2e\c11397f3-0a73-4974-93fe-f54f7fc4bc9d.md"
$wsOverview.Range("C5").Value = ".md"
$wsOverview.Range("D5").Value = "'"
$wsOverview.Range("E5").Value = "Ready for handoff"
$wsOverview.Range("F5").Value = "Ready for handoff"
$wsOverview.Range("G5").Value = "2016-09-01 06:46:14"

# Styles: column B is the hyperlink-styled column, column G is date-styled
$wsOverview.Range("B4").Style = "HyperLink"
$wsOverview.Range("B5").Style = "HyperLink"
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOverview.Range("G5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Hyperlinks
$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/787dd08a-b847-4b32-8a73-c840274ea705.md", "", "", "e2e\787dd08a-b847-4b32-8a73-c840274ea705.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/c11397f3-0a73-4974-93fe-f54f7fc4bc9d.md", "", "", "e2e\c11397f3-0a73-4974-93fe-f54f7fc4bc9d.md") | Out-Null

# -----------------------------------------------------------------------
# Sheet "zh-cn" (sheet2): add two rows
# -----------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$tZh = $wsZh.ListObjects.Item("zh-cn")
$tZh.ListRows.Add() | Out-Null
$tZh.ListRows.Add() | Out-Null

# Row 4
$wsZh.Range("A4").Value = "787dd08a-b847-4b32-8a73-c840274ea705.md"
$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = "Ready for handoff"
$wsZh.Range("D4").Value = "e2e"
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("F4").Value = "'False"
$wsZh.Range("G4").Value = "787dd08a-b847-4b32-8a73-c840274ea705.855a4fa887f5fe69334ba535a5c7693d339f1ba2.zh-cn.xlf"
$wsZh.Range("H4").Value = "2016-09-01 06:46:07"
$wsZh.Range("I4").Value = "'"
$wsZh.Range("J4").Value = "'"
$wsZh.Range("K4").Value = "0001-01-01 00:00:00"
$wsZh.Range("L4").Value = "'"
$wsZh.Range("M4").Value = "'True"
$wsZh.Range("N4").Value = "'"
$wsZh.Range("O4").Value = "'False"
$wsZh.Range("P4").Value = "'"

# Row 5
$wsZh.Range("A5").Value = "c11397f3-0a73-4974-93fe-f54f7fc4bc9d.md"
$wsZh.Range("B5").Value = ".md"
$wsZh.Range("C5").Value = "Ready for handoff"
$wsZh.Range("D5").Value = "e2e"
$wsZh.Range("E5").Value = "ht"
$wsZh.Range("F5").Value = "'False"
$wsZh.Range("G5").Value = "c11397f3-0a73-4974-93fe-f54f7fc4bc9d.5896fc42ba9f5c8baba3c17b17760dc863de37d8.zh-cn.xlf"
$wsZh.Range("H5").Value = "2016-09-01 06:46:07"
$wsZh.Range("I5").Value = "'"
$wsZh.Range("J5").Value = "'"
$wsZh.Range("K5").Value = "0001-01-01 00:00:00"
$wsZh.Range("L5").Value = "'"
$wsZh.Range("M5").Value = "'True"
$wsZh.Range("N5").Value = "'"
$wsZh.Range("O5").Value = "'False"
$wsZh.Range("P5").Value = "'"

# Styles: column A is hyperlink-styled, columns H and K are date-styled
$wsZh.Range("A4").Style = "HyperLink"
$wsZh.Range("A5").Style = "HyperLink"
$wsZh.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("K5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Hyperlinks
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/787dd08a-b847-4b32-8a73-c840274ea705.md", "", "", "787dd08a-b847-4b32-8a73-c840274ea705.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/c11397f3-0a73-4974-93fe-f54f7fc4bc9d.md", "", "", "c11397f3-0a73-4974-93fe-f54f7fc4bc9d.md") | Out-Null

# -----------------------------------------------------------------------
# Sheet "de-de" (sheet3): add two rows
# -----------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$tDe = $wsDe.ListObjects.Item("de-de")
$tDe.ListRows.Add() | Out-Null
$tDe.ListRows.Add() | Out-Null

# Row 4
$wsDe.Range("A4").Value = "787dd08a-b847-4b32-8a73-c840274ea705.md"
$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = "Ready for handoff"
$wsDe.Range("D4").Value = "e2e"
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("F4").Value = "'False"
$wsDe.Range("G4").Value = "787dd08a-b847-4b32-8a73-c840274ea705.855a4fa887f5fe69334ba535a5c7693d339f1ba2.de-de.xlf"
$wsDe.Range("H4").Value = "2016-09-01 06:46:14"
$wsDe.Range("I4").Value = "'"
$wsDe.Range("J4").Value = "'"
$wsDe.Range("K4").Value = "0001-01-01 00:00:00"
$wsDe.Range("L4").Value = "'"
$wsDe.Range("M4").Value = "'True"
$wsDe.Range("N4").Value = "'"
$wsDe.Range("O4").Value = "'False"
$wsDe.Range("P4").Value = "'"

# Row 5
$wsDe.Range("A5").Value = "c11397f3-0a73-4974-93fe-f54f7fc4bc9d.md"
$wsDe.Range("B5").Value = ".md"
$wsDe.Range("C5").Value = "Ready for handoff"
$wsDe.Range("D5").Value = "e2e"
$wsDe.Range("E5").Value = "ht"
$wsDe.Range("F5").Value = "'False"
$wsDe.Range("G5").Value = "c11397f3-0a73-4974-93fe-f54f7fc4bc9d.5896fc42ba9f5c8baba3c17b17760dc863de37d8.de-de.xlf"
$wsDe.Range("H5").Value = "2016-09-01 06:46:14"
$wsDe.Range("I5").Value = "'"
$wsDe.Range("J5").Value = "'"
$wsDe.Range("K5").Value = "0001-01-01 00:00:00"
$wsDe.Range("L5").Value = "'"
$wsDe.Range("M5").Value = "'True"
$wsDe.Range("N5").Value = "'"
$wsDe.Range("O5").Value = "'False"
$wsDe.Range("P5").Value = "'"

# Styles: column A is hyperlink-styled, columns H and K are date-styled
$wsDe.Range("A4").Style = "HyperLink"
$wsDe.Range("A5").Style = "HyperLink"
$wsDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("K5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Hyperlinks
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/787dd08a-b847-4b32-8a73-c840274ea705.md", "", "", "787dd08a-b847-4b32-8a73-c840274ea705.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/c11397f3-0a73-4974-93fe-f54f7fc4bc9d.md", "", "", "c11397f3-0a73-4974-93fe-f54f7fc4bc9d.md") | Out-Null
